$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.191.06"
$ws.Range("E2").Value = "  +1.87%  "

$ws.Range("D3").Value = "2.504.97"
$ws.Range("E3").Value = "  +0.59%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'321.36"

$ws.Range("D6").Value = "'108.25"
$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("D7").Value = "'0.527"
$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.541"
$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("D10").Value = "'39.94"
$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("D11").Value = "'20.24"
$ws.Range("E11").Value = "  +9.26%  "

$ws.Range("D12").Value = "'0.0818"
$ws.Range("E12").Value = "  +1.01%  "

$ws.Range("E13").Value = "  +0.09%  "

$ws.Range("D14").Value = "'7.18"
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("D15").Value = "2.896.02"
$ws.Range("E15").Value = "  +0.57%  "

$ws.Range("D16").Value = "2.506.25"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").Value = "'0.846"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").Value = "48.037.11"
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("D19").Value = "'13.11"
$ws.Range("E19").Value = "  -2.03%  "

$ws.Range("D20").Value = "'6.75"
$ws.Range("E20").Value = "  +1.84%  "

$ws.Range("D21").Value = "0.0₃0947"
$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("E22").Value = "  +1.58%  "

$ws.Range("D23").Value = "'279.39"
$ws.Range("E23").Value = "  +13.15%  "

$ws.Range("D24").Value = "'72.00"
$ws.Range("E24").Value = "  +1.95%  "

$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").Value = "'25.79"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = "  -0.84%  "

$ws.Range("D29").Value = "'9.79"
$ws.Range("E29").Value = "  -1.72%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").Value = "'35.27"
$ws.Range("E31").Value = "  +2.10%  "

$ws.Range("D32").Value = "'49.31"
$ws.Range("E32").Value = "  -1.13%  "

$ws.Range("D33").Value = "'19.53"
$ws.Range("E33").Value = "  -4.48%  "

$ws.Range("D34").Value = "'5.36"
$ws.Range("E34").Value = "  +0.92%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("D37").Value = "'1.96"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").Value = "'4.66"
$ws.Range("E38").Value = "  -2.03%  "

$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").Value = "'0.112"
$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("D41").Value = "'121.35"
$ws.Range("E41").Value = "  +1.77%  "

$ws.Range("E42").Value = "  +0.64%  "

$ws.Range("D43").Value = "'21.55"
$ws.Range("E43").Value = "  -5.10%  "

$ws.Range("D44").Value = "'0.0304"
$ws.Range("E44").Value = "  +2.59%  "

$ws.Range("D45").Value = "2.010.03"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("D46").Value = "'3.18"
$ws.Range("E46").Value = "  +4.91%  "

$ws.Range("E47").Value = "  +4.18%  "

$ws.Range("E48").Value = "  -2.30%  "

$ws.Range("D49").Value = "'8.97"
$ws.Range("E49").Value = "  -1.36%  "

$ws.Range("D50").Value = "'5.18"
$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("D51").Value = "'80.30"
$ws.Range("E51").Value = "  +3.66%  "
